# Efna4-Epha3.xlsx update: re-run of NATMI LR-pair stats with the new TPM matrix.
# "Resolving-Mac" is now included as a sending cluster (in addition to already
# being a target cluster), so every Sending x Target combination over
# {ECs, FAPs, MuSCs, Resolving-Mac} is recomputed for the Efna4 -> Epha3 pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One entry per data row (A2:T13): Sending cluster, Ligand symbol, Receptor
# symbol, Target cluster, then the 16 numeric NATMI statistic columns (E..T).
$rows = @(
    ("ECs","Efna4","Epha3","ECs",3,1,1.180761,3.542283,0.6103536098015441,0.6103536098015441,2,0.6666666666666666,0.114581,0.343743,0.007635610029470834,0.007635610029470834,0.135292776141,1.217634985269,0.004660422144524398,0.004660422144524398),
    ("ECs","Efna4","Epha3","FAPs",3,1,1.180761,3.542283,0.6103536098015441,0.6103536098015441,3,1,14.85444966666667,44.563349,0.9898917347297518,0.9898917347297518,17.539554842863,157.855993585767,0.6041839936050165,0.6041839936050165),
    ("ECs","Efna4","Epha3","MuSCs",3,1,1.180761,3.542283,0.6103536098015441,0.6103536098015441,1,0.3333333333333333,0.037105,0.111315,0.0024726552407774,0.0024726552407774,0.043812136905,0.394309232145,0.001509194052003192,0.001509194052003192),
    ("FAPs","Efna4","Epha3","ECs",3,1,0.4799286666666667,1.439786,0.2480825451952105,0.2480825451952105,2,0.6666666666666666,0.114581,0.343743,0.007635610029470834,0.007635610029470834,0.05499070655533334,0.494916358998,0.001894261570229201,0.001894261570229201),
    ("FAPs","Efna4","Epha3","FAPs",3,1,0.4799286666666667,1.439786,0.2480825451952105,0.2480825451952105,3,1,14.85444966666667,44.563349,0.9898917347297518,0.9898917347297518,7.129076222590445,64.161686003314,0.245574861019459,0.245574861019459),
    ("FAPs","Efna4","Epha3","MuSCs",3,1,0.4799286666666667,1.439786,0.2480825451952105,0.2480825451952105,1,0.3333333333333333,0.037105,0.111315,0.0024726552407774,0.0024726552407774,0.01780775317666667,0.16026977859,0.0006134226055223335,0.0006134226055223335),
    ("MuSCs","Efna4","Epha3","ECs",3,1,0.2047916666666667,0.614375,0.1058599775968842,0.1058599775968842,2,0.6666666666666666,0.114581,0.343743,0.007635610029470834,0.007635610029470834,0.02346523395833333,0.211187105625,0.0008083055066583266,0.0008083055066583265),
    ("MuSCs","Efna4","Epha3","FAPs",3,1,0.2047916666666667,0.614375,0.1058599775968842,0.1058599775968842,3,1,14.85444966666667,44.563349,0.9898917347297518,0.9898917347297518,3.042067504652778,27.378607541875,0.1047899168618323,0.1047899168618323),
    ("MuSCs","Efna4","Epha3","MuSCs",3,1,0.2047916666666667,0.614375,0.1058599775968842,0.1058599775968842,1,0.3333333333333333,0.037105,0.111315,0.0024726552407774,0.0024726552407774,0.007598794791666667,0.068389153125,0.0002617552283935138,0.0002617552283935138),
    ("Resolving-Mac","Efna4","Epha3","ECs",1,0.3333333333333333,0.06907100000000001,0.207213,0.03570386740636119,0.03570386740636119,2,0.6666666666666666,0.114581,0.343743,0.007635610029470834,0.007635610029470834,0.007914224251000001,0.071228018259,0.0002726208080589084,0.0002726208080589084),
    ("Resolving-Mac","Efna4","Epha3","FAPs",1,0.3333333333333333,0.06907100000000001,0.207213,0.03570386740636119,0.03570386740636119,3,1,14.85444966666667,44.563349,0.9898917347297518,0.9898917347297518,1.026011692926333,9.234105236337001,0.03534296324344393,0.03534296324344393),
    ("Resolving-Mac","Efna4","Epha3","MuSCs",1,0.3333333333333333,0.06907100000000001,0.207213,0.03570386740636119,0.03570386740636119,1,0.3333333333333333,0.037105,0.111315,0.0024726552407774,0.0024726552407774,0.002562879455,0.023065915095,0.00008828335485836041,0.00008828335485836041)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowData[$j]
    }
}
